$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 64: qc_flag by record ID for "species_not_normalized"
$ws.Range("A64").Value = "species_not_normalized"
$ws.Range("B64").Value = "Subjects"
$ws.Range("C64").Value = "Species not found (no normalization match)"
$ws.Range("D64").Value = "Soft Stop (Dictionary Update)"

# Update view: scroll so row 47 is at top, select D64 (matches end-user's navigation to new row)
$ws.Range("D64").Select()
